# Report updated to include chunks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the experiment labels in column A:
#  - the old "SequentialOld"/"ParallelOld" rows become the current baseline
#    ("Sequential"/"Parallel")
#  - the previous "Sequential" row becomes the new condensed-chunk result,
#    with an updated timing value, and the previous "Parallel" row becomes
#    "ParallelCondensed"
$ws.Range("A2").Value = "Sequential"
$ws.Range("A3").Value = "SequentialCondensed"
$ws.Range("B3").Value = 16.116607999999999
$ws.Range("A4").Value = "Parallel"
$ws.Range("A5").Value = "ParallelCondensed"

# Column A now needs to fit the longer label text
$ws.Columns.Item(1).ColumnWidth = 20

# Update the remembered selection / active cell for the sheet
$ws.Range("A10").Select()
